$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update random_forest row (row 2)
$ws.Range("B2").Value = 4.8122917338773767
$ws.Range("C2").Value = 0.5101509875805762
$ws.Range("D2").Value = 3.9053775524491159
$ws.Range("E2").Value = 0.47428331762389353
$ws.Range("F2").Value = 0.68868230529315444
$ws.Range("G2").Value = 0.62820896556114525
$ws.Range("H2").Value = 0.52571668237610647
$ws.Range("I2").Value = 0.79049331946619228

# Update lsboost row (row 3)
$ws.Range("B3").Value = 4.884238670532028
$ws.Range("C3").Value = 0.51777808145133675
$ws.Range("D3").Value = 3.8457783387930857
$ws.Range("E3").Value = 0.48857102756581311
$ws.Range("F3").Value = 0.6989785601617643
$ws.Range("G3").Value = 0.61862198969100612
$ws.Range("H3").Value = 0.51142897243418695
$ws.Range("I3").Value = 0.74651559658417987

# Update neural_network row (row 4)
$ws.Range("B4").Value = 4.5424889404325857
$ws.Range("C4").Value = 0.4815491967625124
$ws.Range("D4").Value = 3.6394059764594817
$ws.Range("E4").Value = 0.4225924281126765
$ws.Range("F4").Value = 0.65007109466017365
$ws.Range("G4").Value = 0.5854254634855689
$ws.Range("H4").Value = 0.5774075718873235
$ws.Range("I4").Value = 0.81103739503932837
